$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.342868566513062
$ws.Range("B1").Value = 1.401240110397339
$ws.Range("C1").Value = 1.198653101921082
$ws.Range("D1").Value = 1.26664936542511
$ws.Range("E1").Value = 0.9988322257995605
